# msz - restructuring control processing -> container
#
# The sheet describes a test process ("Record/Dialog" table). A new
# sub-process "003_Profil_003_Abwesenheiten_Normalfall_Anlage" is added as a
# third row, and two new columns are inserted (F:G) to hold the steps that
# belong to the new "popAbwesenheitAnlegen" pop-up dialog.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G - everything that used to live in F:I moves
# to H:K, making room for the new "popAbwesenheitAnlegen" dialog columns.
$ws.Columns("F:G").Insert()

# New row describing the new process. Fill A3/D3/E3 first so that the new
# shared-strings get the same ordering as the authored workbook.
$ws.Range("A3").Value = "003_Profil_003_Abwesenheiten_Normalfall_Anlage"
$ws.Range("D3").Value = "Auswahl Page Abwesenheiten"
$ws.Range("E3").Value = "Abwesenheit anlegen Page Abwesenheiten"

# Header row for the two newly inserted columns - both belong to the new
# "popAbwesenheitAnlegen" dialog.
$ws.Range("F1").Value = "popAbwesenheitAnlegen"
$ws.Range("G1").Value = "popAbwesenheitAnlegen"

# Remaining new-row detail cells.
$ws.Range("F3").Value = "Check defaults"
$ws.Range("G3").Value = "Abwesenheit Sylvester anlegen"
$ws.Range("H3").Value = "Abwesenheit Sylvester checken"
$ws.Range("I3").Value = "Abwesenheit Sylvester löschen"
$ws.Range("K3").Value = "Menueauswahl Überblick"

# G3 carries an explicit Text number format in the authored workbook.
$ws.Range("G3").NumberFormat = "@"

# Resize column E to fit its (now longer) content, then mirror that width
# onto the two freshly inserted columns F and G.
$ws.Columns("E").AutoFit()
$newWidth = $ws.Columns("E").ColumnWidth
$ws.Columns("F").ColumnWidth = $newWidth
$ws.Columns("G").ColumnWidth = $newWidth

# Page setup was switched to portrait / paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to the next empty row, like a user would after
# finishing data entry.
$ws.Range("A4").Select() | Out-Null
